$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so values like "43.20" or "0.0000111"
# are stored verbatim instead of being reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.300.52'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.229.45'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '244.96'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '0.628'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').Value = '74.41'
$ws.Range('E7').Value = '  -2.26%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = '43.20'
$ws.Range('E10').Value = '  +4.44%  '
$ws.Range('D11').Value = '0.0961'
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '14.37'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').Value = '0.850'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '2.234.29'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '42.141.49'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').Value = '0.0000111'
$ws.Range('E18').Value = '  +13.67%  '
$ws.Range('D19').Value = '6.17'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('D21').Value = '10.14'
$ws.Range('E21').Value = '  +38.25%  '
$ws.Range('D22').Value = '231.11'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').Value = '  -4.71%  '
$ws.Range('D24').Value = '11.76'
$ws.Range('E24').Value = '  +5.30%  '
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').Value = '3.65'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D29').Value = '167.12'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').Value = '5.84'
$ws.Range('E31').Value = '  +18.51%  '
$ws.Range('D32').Value = '0.0804'
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '0.125'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '29.60'
$ws.Range('E35').Value = '  -7.75%  '
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').Value = '0.0306'
$ws.Range('E37').Value = '  +2.62%  '
$ws.Range('D38').Value = '13.20'
$ws.Range('E38').Value = '  -4.98%  '
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = '5.62'
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('D41').Value = '63.28'
$ws.Range('E41').Value = '  +4.84%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('E43').Value = '  +2.30%  '
$ws.Range('D44').Value = '104.81'
$ws.Range('E44').Value = '  -6.34%  '
$ws.Range('E45').Value = '  +3.59%  '
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('E47').Value = '  +6.34%  '
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('D49').Value = '1.17'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('E51').Value = '  -1.32%  '
